# Update column F (dSF) values on Sheet1 to reflect repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -4
    7  = -9
    9  = 11
    10 = 0
    11 = -2
    12 = 9
    13 = -4
    23 = -6
    24 = -7
    25 = -3
    27 = -1
    29 = -1
    31 = 1
    32 = -3
    35 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
